$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 27.09457833333333
$ws.Range("H2").Value = 81.28373500000001
$ws.Range("I2").Value = 0.9716243152921379
$ws.Range("J2").Value = 0.9716243152921381
$ws.Range("M2").Value = 105.0653766666667
$ws.Range("N2").Value = 315.19613
$ws.Range("O2").Value = 0.9818992286940446
$ws.Range("P2").Value = 0.9818992286940446
$ws.Range("Q2").Value = 2846.702078216172
$ws.Range("R2").Value = 25620.31870394555
$ws.Range("S2").Value = 0.9540371657657294
$ws.Range("T2").Value = 0.9540371657657296
$ws.Range("G3").Value = 27.09457833333333
$ws.Range("H3").Value = 81.28373500000001
$ws.Range("I3").Value = 0.9716243152921379
$ws.Range("J3").Value = 0.9716243152921381
$ws.Range("O3").Value = 0.0115004770447132
$ws.Range("P3").Value = 0.0115004770447132
$ws.Range("Q3").Value = 33.34194685864611
$ws.Range("R3").Value = 300.0775217278151
$ws.Range("S3").Value = 0.01117414313410241
$ws.Range("T3").Value = 0.01117414313410241
$ws.Range("G4").Value = 27.09457833333333
$ws.Range("H4").Value = 81.28373500000001
$ws.Range("I4").Value = 0.9716243152921379
$ws.Range("J4").Value = 0.9716243152921381
$ws.Range("M4").Value = 0.3093633333333333
$ws.Range("N4").Value = 0.92809
$ws.Range("O4").Value = 0.002891186687979499
$ws.Range("P4").Value = 0.0028911866879795
$ws.Range("Q4").Value = 8.382069068461112
$ws.Range("R4").Value = 75.43862161615
$ws.Range("S4").Value = 0.002809147286089825
$ws.Range("T4").Value = 0.002809147286089826
$ws.Range("G5").Value = 27.09457833333333
$ws.Range("H5").Value = 81.28373500000001
$ws.Range("I5").Value = 0.9716243152921379
$ws.Range("J5").Value = 0.9716243152921381
$ws.Range("M5").Value = 0.3968826666666667
$ws.Range("N5").Value = 1.190648
$ws.Range("O5").Value = 0.003709107573262739
$ws.Range("P5").Value = 0.003709107573262739
$ws.Range("Q5").Value = 10.75336850114222
$ws.Range("R5").Value = 96.78031651028
$ws.Range("S5").Value = 0.003603859106216291
$ws.Range("T5").Value = 0.003603859106216292
$ws.Range("I6").Value = 0.01492035819923828
$ws.Range("J6").Value = 0.01492035819923829
$ws.Range("M6").Value = 105.0653766666667
$ws.Range("N6").Value = 315.19613
$ws.Range("O6").Value = 0.9818992286940446
$ws.Range("P6").Value = 0.9818992286940446
$ws.Range("Q6").Value = 43.71423607357
$ws.Range("R6").Value = 393.42812466213
$ws.Range("S6").Value = 0.01465028820767094
$ws.Range("T6").Value = 0.01465028820767094
$ws.Range("I7").Value = 0.01492035819923828
$ws.Range("J7").Value = 0.01492035819923829
$ws.Range("O7").Value = 0.0115004770447132
$ws.Range("P7").Value = 0.0115004770447132
$ws.Range("R7").Value = 4.608019829529001
$ws.Range("S7").Value = 0.0001715912369692382
$ws.Range("T7").Value = 0.0001715912369692382
$ws.Range("I8").Value = 0.01492035819923828
$ws.Range("J8").Value = 0.01492035819923829
$ws.Range("M8").Value = 0.3093633333333333
$ws.Range("N8").Value = 0.92809
$ws.Range("O8").Value = 0.002891186687979499
$ws.Range("P8").Value = 0.0028911866879795
$ws.Range("Q8").Value = 0.12871587401
$ws.Range("R8").Value = 1.15844286609
$ws.Range("S8").Value = 0.0000431375410055235
$ws.Range("T8").Value = 0.00004313754100552351
$ws.Range("I9").Value = 0.01492035819923828
$ws.Range("J9").Value = 0.01492035819923829
$ws.Range("M9").Value = 0.3968826666666667
$ws.Range("N9").Value = 1.190648
$ws.Range("O9").Value = 0.003709107573262739
$ws.Range("P9").Value = 0.003709107573262739
$ws.Range("Q9").Value = 0.165129780472
$ws.Range("R9").Value = 1.486168024248
$ws.Range("S9").Value = 0.00005534121359258752
$ws.Range("T9").Value = 0.00005534121359258752
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.3752133333333333
$ws.Range("H10").Value = 1.12564
$ws.Range("I10").Value = 0.01345532650862368
$ws.Range("J10").Value = 0.01345532650862368
$ws.Range("M10").Value = 105.0653766666667
$ws.Range("N10").Value = 315.19613
$ws.Range("O10").Value = 0.9818992286940446
$ws.Range("P10").Value = 0.9818992286940446
$ws.Range("Q10").Value = 39.42193019702222
$ws.Range("R10").Value = 354.7973717732
$ws.Range("S10").Value = 0.01321177472064412
$ws.Range("T10").Value = 0.01321177472064412
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.3752133333333333
$ws.Range("H11").Value = 1.12564
$ws.Range("I11").Value = 0.01345532650862368
$ws.Range("J11").Value = 0.01345532650862368
$ws.Range("O11").Value = 0.0115004770447132
$ws.Range("P11").Value = 0.0115004770447132
$ws.Range("Q11").Value = 0.4617286479511111
$ws.Range("R11").Value = 4.15555783156
$ws.Range("S11").Value = 0.0001547426736415475
$ws.Range("T11").Value = 0.0001547426736415475
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.3752133333333333
$ws.Range("H12").Value = 1.12564
$ws.Range("I12").Value = 0.01345532650862368
$ws.Range("J12").Value = 0.01345532650862368
$ws.Range("M12").Value = 0.3093633333333333
$ws.Range("N12").Value = 0.92809
$ws.Range("O12").Value = 0.002891186687979499
$ws.Range("P12").Value = 0.0028911866879795
$ws.Range("Q12").Value = 0.1160772475111111
$ws.Range("R12").Value = 1.0446952276
$ws.Range("S12").Value = 0.00003890186088415044
$ws.Range("T12").Value = 0.00003890186088415045
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.3752133333333333
$ws.Range("H13").Value = 1.12564
$ws.Range("I13").Value = 0.01345532650862368
$ws.Range("J13").Value = 0.01345532650862368
$ws.Range("M13").Value = 0.3968826666666667
$ws.Range("N13").Value = 1.190648
$ws.Range("O13").Value = 0.003709107573262739
$ws.Range("P13").Value = 0.003709107573262739
$ws.Range("Q13").Value = 0.1489156683022222
$ws.Range("R13").Value = 1.34024101472
$ws.Range("S13").Value = 0.00004990725345385896
$ws.Range("T13").Value = 0.00004990725345385896
